$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Mark every remaining test case (rows 4-26) to run by setting Runmode
# column (D) to "Y" instead of "N" -- running all notification test cases.
$ws.Range("D4:D26").Value = "Y"

# Keep the selection in sync with the edited range, matching what Excel
# records after selecting D3:D26 and editing it.
$ws.Range("D3:D26").Select()
